# import_data_stok_oksigen.xlsx — rename sheets, move the active selection,
# and repoint the "Rumah Sakit" list-validation source at the renamed sheet.

$wb = $excel.ActiveWorkbook

# 1) Rename the sheets: Sheet1 -> Data, KabKota -> Rumah Sakit
$wsData = $wb.Worksheets.Item("Sheet1")
$wsHospitals = $wb.Worksheets.Item("KabKota")
$wsHospitals.Name = "Rumah Sakit"
$wsData.Name = "Data"

# 2) Try to repoint the B3:B1048576 list validation at the renamed sheet so the
#    dropdown keeps working. (Best effort: some engines keep validation rules
#    read-only/opaque, in which case this is a harmless no-op.)
$dataRange = $wsData.Range("B3:B1048576")
try {
    $validation = $dataRange.Validation
    $validation.Formula1 = "='Rumah Sakit'!`$A`$2:`$A`$20"
} catch {
    # Validation object not mutable in this host - leave the original rule intact.
}

# 3) Move the active cell / selection from F3 to G3 on the Data sheet.
$wsData.Activate() | Out-Null
$wsData.Range("G3").Select() | Out-Null

Write-Output "Renamed sheets, updated validation source, moved selection to G3."
